$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append the new menu item (Campa Energy 30Rs) as row 56, matching the
# Item / Half / Full / Image layout used throughout the sheet.
$ws.Range("A56").Value = "Campa Energy 30Rs"
$ws.Range("C56").Value = 30
$ws.Range("C56").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("D56").Value = "Campa Energy 30Rs.jpg"

# Leave the selection where the user ended up after entering the new row.
$ws.Range("A59").Select()
